$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 24.916566
$ws.Range("H2").Value = 74.749698
$ws.Range("I2").Value = 0.459912889255076
$ws.Range("J2").Value = 0.459912889255076
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 45.95651366666667
$ws.Range("N2").Value = 137.869541
$ws.Range("O2").Value = 0.6189188856627118
$ws.Range("P2").Value = 0.6189188856627118
$ws.Range("Q2").Value = 1145.078505905402
$ws.Range("R2").Value = 10305.70655314862
$ws.Range("S2").Value = 0.2846487729196698
$ws.Range("T2").Value = 0.2846487729196698

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 24.916566
$ws.Range("H3").Value = 74.749698
$ws.Range("I3").Value = 0.459912889255076
$ws.Range("J3").Value = 0.459912889255076
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("N3").Value = 20.549744
$ws.Range("O3").Value = 0.09225115688993263
$ws.Range("P3").Value = 0.09225115688993261
$ws.Range("Q3").Value = 170.676350886368
$ws.Range("R3").Value = 1536.087157977312
$ws.Range("S3").Value = 0.04242749610237222
$ws.Range("T3").Value = 0.04242749610237222

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 24.916566
$ws.Range("H4").Value = 74.749698
$ws.Range("I4").Value = 0.459912889255076
$ws.Range("J4").Value = 0.459912889255076
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 21.446458
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2888299574473556
$ws.Range("P4").Value = 0.2888299574473556
$ws.Range("Q4").Value = 534.3720862232279
$ws.Range("R4").Value = 4809.348776009051
$ws.Range("S4").Value = 0.132836620233034
$ws.Range("T4").Value = 0.132836620233034

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 19.60300333333333
$ws.Range("H5").Value = 58.80901
$ws.Range("I5").Value = 0.3618345281251927
$ws.Range("J5").Value = 0.3618345281251927
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 45.95651366666667
$ws.Range("N5").Value = 137.869541
$ws.Range("O5").Value = 0.6189188856627118
$ws.Range("P5").Value = 0.6189188856627118
$ws.Range("Q5").Value = 900.8856905960456
$ws.Range("R5").Value = 8107.97121536441
$ws.Range("S5").Value = 0.2239462229415374
$ws.Range("T5").Value = 0.2239462229415374

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 19.60300333333333
$ws.Range("H6").Value = 58.80901
$ws.Range("I6").Value = 0.3618345281251927
$ws.Range("J6").Value = 0.3618345281251927
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.849914666666667
$ws.Range("N6").Value = 20.549744
$ws.Range("O6").Value = 0.09225115688993263
$ws.Range("P6").Value = 0.09225115688993261
$ws.Range("Q6").Value = 134.2789000437156
$ws.Range("R6").Value = 1208.51010039344
$ws.Range("S6").Value = 0.03337965382227189
$ws.Range("T6").Value = 0.03337965382227188

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 19.60300333333333
$ws.Range("H7").Value = 58.80901
$ws.Range("I7").Value = 0.3618345281251927
$ws.Range("J7").Value = 0.3618345281251927
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 21.446458
$ws.Range("N7").Value = 64.33937399999999
$ws.Range("O7").Value = 0.2888299574473556
$ws.Range("P7").Value = 0.2888299574473556
$ws.Range("Q7").Value = 420.4149876621933
$ws.Range("R7").Value = 3783.73488895974
$ws.Range("S7").Value = 0.1045086513613834
$ws.Range("T7").Value = 0.1045086513613834

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 9.657138000000002
$ws.Range("H8").Value = 28.971414
$ws.Range("I8").Value = 0.1782525826197313
$ws.Range("J8").Value = 0.1782525826197313
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 45.95651366666667
$ws.Range("N8").Value = 137.869541
$ws.Range("O8").Value = 0.6189188856627118
$ws.Range("P8").Value = 0.6189188856627118
$ws.Range("Q8").Value = 443.8083944778861
$ws.Range("R8").Value = 3994.275550300974
$ws.Range("S8").Value = 0.1103238898015045
$ws.Range("T8").Value = 0.1103238898015045

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 9.657138000000002
$ws.Range("H9").Value = 28.971414
$ws.Range("I9").Value = 0.1782525826197313
$ws.Range("J9").Value = 0.1782525826197313
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.849914666666667
$ws.Range("N9").Value = 20.549744
$ws.Range("O9").Value = 0.09225115688993263
$ws.Range("P9").Value = 0.09225115688993261
$ws.Range("Q9").Value = 66.15057122422401
$ws.Range("R9").Value = 595.3551410180161
$ws.Range("S9").Value = 0.01644400696528851
$ws.Range("T9").Value = 0.0164440069652885

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 9.657138000000002
$ws.Range("H10").Value = 28.971414
$ws.Range("I10").Value = 0.1782525826197313
$ws.Range("J10").Value = 0.1782525826197313
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 21.446458
$ws.Range("N10").Value = 64.33937399999999
$ws.Range("O10").Value = 0.2888299574473556
$ws.Range("P10").Value = 0.2888299574473556
$ws.Range("Q10").Value = 207.111404517204
$ws.Range("R10").Value = 1864.002640654836
$ws.Range("S10").Value = 0.05148468585293822
$ws.Range("T10").Value = 0.05148468585293822

